$d = $word.ActiveDocument

# --- Edit 1: "virtualbox" -> "VirtualBox" in the chocolatey-pin sentence ---
$sel = $d.Range(0, 0)
$found = $sel.Find.Execute("virtualbox version", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $wordStart = $sel.Start
    $wordEnd = $wordStart + 10
    $r = $d.Range($wordStart, $wordEnd)
    if ($r.Text -ceq "virtualbox") {
        $r.Text = "VirtualBox"
    }
}

# --- Edit 2: insert "~" right after "ubuntu-vanilla. " (merging the "." and " " runs into ". ~") ---
$sel = $d.Range(0, 0)
$found = $sel.Find.Execute("ubuntu-vanilla", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $afterName = $sel.End
    $r2 = $d.Range($afterName, $afterName + 2)
    if ($r2.Text -ceq ". ") {
        $r2.Text = ". ~"
    }
}

# --- Edit 3: append "~" after "command:" at the end of the same paragraph ---
$sel = $d.Range(0, 0)
$found = $sel.Find.Execute("vagrant init command:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $cmdStart = $sel.End - 8
    $cmdEnd = $sel.End
    $r3 = $d.Range($cmdStart, $cmdEnd)
    if ($r3.Text -ceq "command:") {
        $r3.InsertAfter("~")
    }
}
